# Apply the "cryptos" price/volume refresh described in the commit:
#   "Updated cryptos list on Mon Apr  3 07:42:40 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 36/37: FraxShare and VeChain swapped rank order with refreshed data ---
$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'

# --- Column D (Price): these are plain-text cells in the source sheet. Whenever the
# new price string would otherwise be auto-parsed as a number by Excel, briefly mark
# the cell as Text so the literal string is preserved exactly, then restore the default
# (Normal) cell style so formatting still matches the rest of the unstyled column.
$ws.Range("D2").Value = '27.902.33'
$ws.Range("D3").Value = '1.792.15'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9999'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5126'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3925'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07820'
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '40.89'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.246'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.000'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.22'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.242'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = '1.782.29'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001077'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06529'
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.931'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = '27.979.73'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.66'
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Value = '1.994.39'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.371'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.25'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1079'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.042'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.617'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.493'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07090'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '8.833'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02304'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2132'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.54'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.017'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6102'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9996'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.152'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.21'
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5921'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.701'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.70'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.208'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.912'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06821'
$ws.Range("D51").Style = "Normal"

# --- Column E (Volume 1h): already text (leading/trailing spaces + %), plain assignment ---
$ws.Range("E2").Value = '  -2.21%  '
$ws.Range("E3").Value = '  -1.89%  '
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("E5").Value = '  -2.00%  '
$ws.Range("E6").Value = '  -0.26%  '
$ws.Range("E7").Value = '  -0.70%  '
$ws.Range("E8").Value = '  +0.91%  '
$ws.Range("E9").Value = '  -7.19%  '
$ws.Range("E11").Value = '  -2.59%  '
$ws.Range("E12").Value = '  -2.91%  '
$ws.Range("E13").Value = '  -0.27%  '
$ws.Range("E14").Value = '  -5.00%  '
$ws.Range("E15").Value = '  -3.89%  '
$ws.Range("E16").Value = '  -2.50%  '
$ws.Range("E17").Value = '  -2.73%  '
$ws.Range("E18").Value = '  -5.12%  '
$ws.Range("E19").Value = '  -1.47%  '
$ws.Range("E20").Value = '  -0.29%  '
$ws.Range("E21").Value = '  -4.06%  '
$ws.Range("E22").Value = '  -2.61%  '
$ws.Range("E23").Value = '  -2.10%  '
$ws.Range("E24").Value = '  -3.60%  '
$ws.Range("E25").Value = '  -2.24%  '
$ws.Range("E26").Value = '  +0.51%  '
$ws.Range("E27").Value = '  -4.50%  '
$ws.Range("E28").Value = '  -2.06%  '
$ws.Range("E29").Value = '  -1.59%  '
$ws.Range("E30").Value = '  +1.08%  '
$ws.Range("E31").Value = '  -1.72%  '
$ws.Range("E32").Value = '  -5.21%  '
$ws.Range("E33").Value = '  -1.51%  '
$ws.Range("E34").Value = '  -4.36%  '
$ws.Range("E35").Value = '  -7.49%  '
$ws.Range("E36").Value = '  +0.98%  '
$ws.Range("E37").Value = '  -3.80%  '
$ws.Range("E38").Value = '  -4.42%  '
$ws.Range("E39").Value = '  +0.45%  '
$ws.Range("E41").Value = '  -4.33%  '
$ws.Range("E42").Value = '  -0.26%  '
$ws.Range("E43").Value = '  -3.14%  '
$ws.Range("E44").Value = '  -2.73%  '
$ws.Range("E45").Value = '  -6.42%  '
$ws.Range("E46").Value = '  -2.21%  '
$ws.Range("E47").Value = '  -2.10%  '
$ws.Range("E48").Value = '  -2.07%  '
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("E50").Value = '  -4.49%  '
$ws.Range("E51").Value = '  -2.40%  '
